$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.553.91'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '1.904.97'
$ws.Range('E3').Value = '  +3.10%  '
$ws.Range('E4').Value = '  +0.72%  '
$ws.Range('D5').Value = "'246.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.66%  '
$ws.Range('D6').Value = "'0.632"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('D8').Value = "'42.22"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +2.92%  '
$ws.Range('E10').Value = '  +1.36%  '
$ws.Range('D11').Value = "'0.0999"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').Value = '2.179.16'
$ws.Range('D13').Value = "'12.42"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.17%  '
$ws.Range('D14').Value = '1.926.82'
$ws.Range('E14').Value = '  +4.35%  '
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = "'4.86"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = '35.560.80'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').Value = "'71.99"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '0.0₃0813'
$ws.Range('E19').Value = '  +2.45%  '
$ws.Range('D20').Value = "'244.00"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.28%  '
$ws.Range('D21').Value = "'12.49"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('E22').Value = '  +2.33%  '
$ws.Range('E23').Value = '  +0.62%  '
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'172.16"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = "'2.21"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +27.34%  '
$ws.Range('D27').Value = "'8.59"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.15%  '
$ws.Range('E28').Value = '  +2.38%  '
$ws.Range('E29').Value = '  +0.61%  '
$ws.Range('D30').Value = "'0.994"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +29.82%  '
$ws.Range('E31').Value = '  +3.40%  '
$ws.Range('D32').Value = "'0.0567"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.84%  '
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').Value = "'4.16"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.38%  '
$ws.Range('E35').Value = '  +6.35%  '
$ws.Range('E36').Value = '  +3.39%  '
$ws.Range('D37').Value = "'1.32"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.74%  '
$ws.Range('E38').Value = '  +2.65%  '
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('D40').Value = "'91.11"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').Value = "'51.01"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +47.84%  '
$ws.Range('D42').Value = '1.353.55'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('D43').Value = "'15.49"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.59%  '
$ws.Range('E44').Value = '  +11.68%  '
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').Value = "'12.60"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = "'6.67"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.73%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = "'2.76"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').Value = '2.090.62'
$ws.Range('E50').Value = '  +3.01%  '
$ws.Range('E51').Value = '  +2.40%  '
